$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: Iy [m^4] -> Iz [m^4]
$ws.Range("D6").Value = "Iz [m^4]"
$ws.Range("I6").Value = "Iz [m^4]"
$ws.Range("N6").Value = "Iz [m^4]"

# Update cross-section numeric values
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 88.73
$ws.Range("E8").Value = 11.06
$ws.Range("F8").Value = 0.44
$ws.Range("H8").Value = 8
$ws.Range("I8").Value = 88.73
$ws.Range("J8").Value = 11.06
$ws.Range("K8").Value = 0.44
$ws.Range("M8").Value = 8
$ws.Range("N8").Value = 88.73
$ws.Range("O8").Value = 11.06
$ws.Range("P8").Value = 0.44
$ws.Range("B9").Value = 13.75
$ws.Range("C9").Value = 7.42
$ws.Range("D9").Value = 70.98
$ws.Range("E9").Value = 10.26
$ws.Range("F9").Value = 0.44
$ws.Range("G9").Value = 13.75
$ws.Range("H9").Value = 7.3
$ws.Range("I9").Value = 67.45999999999999
$ws.Range("J9").Value = 10.09
$ws.Range("K9").Value = 0.44
$ws.Range("L9").Value = 13.75
$ws.Range("M9").Value = 7.3
$ws.Range("N9").Value = 67.45999999999999
$ws.Range("O9").Value = 10.09
$ws.Range("P9").Value = 0.44
$ws.Range("B10").Value = 27.5
$ws.Range("C10").Value = 6.85
$ws.Range("D10").Value = 55.77
$ws.Range("E10").Value = 9.470000000000001
$ws.Range("F10").Value = 0.44
$ws.Range("G10").Value = 27.5
$ws.Range("H10").Value = 6.6
$ws.Range("I10").Value = 49.9
$ws.Range("J10").Value = 9.119999999999999
$ws.Range("K10").Value = 0.44
$ws.Range("L10").Value = 27.5
$ws.Range("M10").Value = 6.6
$ws.Range("N10").Value = 49.9
$ws.Range("O10").Value = 9.119999999999999
$ws.Range("P10").Value = 0.44
$ws.Range("B11").Value = 41.25
$ws.Range("C11").Value = 6.28
$ws.Range("D11").Value = 42.9
$ws.Range("E11").Value = 8.67
$ws.Range("F11").Value = 0.44
$ws.Range("G11").Value = 41.25
$ws.Range("H11").Value = 5.9
$ws.Range("I11").Value = 35.68
$ws.Range("J11").Value = 8.16
$ws.Range("K11").Value = 0.44
$ws.Range("L11").Value = 41.25
$ws.Range("M11").Value = 5.9
$ws.Range("N11").Value = 35.68
$ws.Range("O11").Value = 8.16
$ws.Range("P11").Value = 0.44
$ws.Range("B12").Value = 55
$ws.Range("C12").Value = 5.7
$ws.Range("D12").Value = 32.19
$ws.Range("E12").Value = 7.88
$ws.Range("F12").Value = 0.44
$ws.Range("G12").Value = 55
$ws.Range("H12").Value = 5.2
$ws.Range("I12").Value = 24.47
$ws.Range("J12").Value = 7.19
$ws.Range("K12").Value = 0.44
$ws.Range("L12").Value = 55
$ws.Range("M12").Value = 5.2
$ws.Range("N12").Value = 24.47
$ws.Range("O12").Value = 7.19
$ws.Range("P12").Value = 0.44
$ws.Range("B13").Value = 68.75
$ws.Range("C13").Value = 5.12
$ws.Range("D13").Value = 23.43
$ws.Range("E13").Value = 7.08
$ws.Range("F13").Value = 0.44
$ws.Range("G13").Value = 68.75
$ws.Range("H13").Value = 4.5
$ws.Range("I13").Value = 15.9
$ws.Range("J13").Value = 6.22
$ws.Range("K13").Value = 0.44
$ws.Range("L13").Value = 68.75
$ws.Range("M13").Value = 4.5
$ws.Range("N13").Value = 15.9
$ws.Range("O13").Value = 6.22
$ws.Range("P13").Value = 0.44
$ws.Range("B14").Value = 82.5
$ws.Range("C14").Value = 4.55
$ws.Range("D14").Value = 16.43
$ws.Range("E14").Value = 6.29
$ws.Range("F14").Value = 0.44
$ws.Range("G14").Value = 82.5
$ws.Range("H14").Value = 4.13
$ws.Range("I14").Value = 12.34
$ws.Range("J14").Value = 5.71
$ws.Range("K14").Value = 0.44
$ws.Range("L14").Value = 82.5
$ws.Range("M14").Value = 4.13
$ws.Range("N14").Value = 12.34
$ws.Range("O14").Value = 5.71
$ws.Range("P14").Value = 0.44
$ws.Range("B15").Value = 96.25
$ws.Range("C15").Value = 3.98
$ws.Range("D15").Value = 10.99
$ws.Range("E15").Value = 5.49
$ws.Range("F15").Value = 0.44
$ws.Range("G15").Value = 96.25
$ws.Range("H15").Value = 3.77
$ws.Range("I15").Value = 9.359999999999999
$ws.Range("J15").Value = 5.21
$ws.Range("K15").Value = 0.44
$ws.Range("L15").Value = 96.25
$ws.Range("M15").Value = 3.77
$ws.Range("N15").Value = 9.359999999999999
$ws.Range("O15").Value = 5.21
$ws.Range("P15").Value = 0.44
$ws.Range("B16").Value = 110
$ws.Range("D16").Value = 6.9
$ws.Range("E16").Value = 4.7
$ws.Range("F16").Value = 0.44
$ws.Range("G16").Value = 110
$ws.Range("I16").Value = 6.9
$ws.Range("J16").Value = 4.7
$ws.Range("K16").Value = 0.44
$ws.Range("L16").Value = 110
$ws.Range("N16").Value = 6.9
$ws.Range("O16").Value = 4.7
$ws.Range("P16").Value = 0.44
